# Generate Report for Handback
# Marks the two handed-back localization files as in sync with en-US, records
# the handback target/file columns (Latest Target File / Latest Handback File)
# for the zh-cn and de-de sheets, and stamps the new handback datetime.

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1172bd143dc9d16488cc4db310f8b5c932f350cf/e2e/7632a710-7701-416b-8d2e-59a3ff7ec853.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1172bd143dc9d16488cc4db310f8b5c932f350cf/e2e/83844d55-54b6-48ac-b901-a46e4587ae1a.md"

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: status text + column widths -----------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn sheet ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

# Latest Handback DateTime column keeps referencing the same timestamp cells,
# just with the updated text.
$zh.Range("K2").Value = "2016-08-21 06:37:43"
$zh.Range("K3").Value = "2016-08-21 06:37:43"

# Latest Target File / Latest Handback File columns for the two rows.
$zh.Hyperlinks.Add($zh.Range("I2"), $urlA, "", "", "7632a710-7701-416b-8d2e-59a3ff7ec853.md")
$zh.Range("J2").Value = "7632a710-7701-416b-8d2e-59a3ff7ec853.f887dead448423eb5b42715e9a8f1e2e06d2c9ed.zh-cn.xlf"

$zh.Hyperlinks.Add($zh.Range("I3"), $urlB, "", "", "83844d55-54b6-48ac-b901-a46e4587ae1a.md")
$zh.Range("J3").Value = "83844d55-54b6-48ac-b901-a46e4587ae1a.a7a24bc970f2a3d0068b61ae0e74343b0d7b1a14.zh-cn.xlf"

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

# Both rows were handed back together, at the new timestamp.
$de.Range("K2").Value = "2016-08-21 06:37:49"
$de.Range("K3").Value = "2016-08-21 06:37:49"

$de.Hyperlinks.Add($de.Range("I2"), $urlA, "", "", "7632a710-7701-416b-8d2e-59a3ff7ec853.md")
$de.Range("J2").Value = "7632a710-7701-416b-8d2e-59a3ff7ec853.f887dead448423eb5b42715e9a8f1e2e06d2c9ed.de-de.xlf"

$de.Hyperlinks.Add($de.Range("I3"), $urlB, "", "", "83844d55-54b6-48ac-b901-a46e4587ae1a.md")
$de.Range("J3").Value = "83844d55-54b6-48ac-b901-a46e4587ae1a.a7a24bc970f2a3d0068b61ae0e74343b0d7b1a14.de-de.xlf"

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
